# Correção nos dados e início da análise PNAD 2009
#
# The "grandes regiões e unidades da federação" header row (row 6) had no
# data of its own; the data table actually starts one row below what it
# should. Removing this row and shifting everything up by one realigns
# each region label with its correct values (and drops the now-unused
# trailing row that falls off the bottom of the table).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A6:G6").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftUp)
